# Atualizado por script em 05-11-2023 08:45
#
# Re-scrapes a handful of Serie A 2023-2024 fixtures: a handful of rows had
# their home/away data shuffled between neighbouring rows (matches recorded
# out of kickoff-time order get re-sorted), and three brand-new matches
# (rows 103-105) are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: snapshot the F:V ("match data") cells of a row into a hashtable
# so we can permute whole rows without clobbering data we still need to
# read from a cell another assignment already overwrote.
# ---------------------------------------------------------------------
function Get-RowSnapshot($row) {
    $snap = @{}
    for ($col = 6; $col -le 22; $col++) {
        $snap[$col] = $ws.Cells.Item($row, $col).Value()
    }
    return $snap
}

function Set-RowSnapshot($row, $snap) {
    for ($col = 6; $col -le 22; $col++) {
        $ws.Cells.Item($row, $col).Value = $snap[$col]
    }
}

# ---------------------------------------------------------------------
# Rows 4 <-> 5: Inter-Monza / Genoa-Fiorentina swapped places
# ---------------------------------------------------------------------
$snap4 = Get-RowSnapshot 4
$snap5 = Get-RowSnapshot 5
Set-RowSnapshot 4 $snap5
Set-RowSnapshot 5 $snap4

# ---------------------------------------------------------------------
# Rows 26 <-> 27: Atalanta-Monza / Napoli-Lazio swapped places
# ---------------------------------------------------------------------
$snap26 = Get-RowSnapshot 26
$snap27 = Get-RowSnapshot 27
Set-RowSnapshot 26 $snap27
Set-RowSnapshot 27 $snap26

# ---------------------------------------------------------------------
# Rows 28 <-> 29: Inter-Fiorentina / Torino-Genoa swapped places
# ---------------------------------------------------------------------
$snap28 = Get-RowSnapshot 28
$snap29 = Get-RowSnapshot 29
Set-RowSnapshot 28 $snap29
Set-RowSnapshot 29 $snap28

# ---------------------------------------------------------------------
# Rows 53 -> 54 -> 55 -> 53 (three-way rotation):
#   new 53 = old 55 (Verona-Atalanta)
#   new 54 = old 53 (Cagliari-AC Milan)
#   new 55 = old 54 (Empoli-Salernitana)
# ---------------------------------------------------------------------
$snap53 = Get-RowSnapshot 53
$snap54 = Get-RowSnapshot 54
$snap55 = Get-RowSnapshot 55
Set-RowSnapshot 53 $snap55
Set-RowSnapshot 54 $snap53
Set-RowSnapshot 55 $snap54

# ---------------------------------------------------------------------
# Rows 86 <-> 87: Bologna-Frosinone / Salernitana-Cagliari swapped places
# ---------------------------------------------------------------------
$snap86 = Get-RowSnapshot 86
$snap87 = Get-RowSnapshot 87
Set-RowSnapshot 86 $snap87
Set-RowSnapshot 87 $snap86

# ---------------------------------------------------------------------
# Append three new fixtures as rows 103-105. Copy row 102's formatting
# (bold/bordered index column, date-formatted kickoff column) down first,
# then fill in the values.
# ---------------------------------------------------------------------
$ws.Range("A102:V102").Copy()
$ws.Range("A103:V105").PasteSpecial(-4122)

$ws.Cells.Item(103, 1).Value = 102
$ws.Cells.Item(103, 2).Value = "italy"
$ws.Cells.Item(103, 3).Value = "serie-a"
$ws.Cells.Item(103, 4).Value = "2023-2024"
$ws.Cells.Item(103, 5).Value = 45234.625
$ws.Cells.Item(103, 6).Value = "Salernitana"
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = "Napoli"
$ws.Cells.Item(103, 9).Value = 2
$ws.Cells.Item(103, 10).Value = 5.76
$ws.Cells.Item(103, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(103, 12).Value = 8.15
$ws.Cells.Item(103, 13).Value = "04/11/2023 14:59"
$ws.Cells.Item(103, 14).Value = 4.37
$ws.Cells.Item(103, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(103, 16).Value = 5.28
$ws.Cells.Item(103, 17).Value = "04/11/2023 14:58"
$ws.Cells.Item(103, 18).Value = 1.57
$ws.Cells.Item(103, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(103, 20).Value = 1.39
$ws.Cells.Item(103, 21).Value = "04/11/2023 14:49"
$ws.Cells.Item(103, 22).Value = "https://www.betexplorer.com/football/italy/serie-a/salernitana-napoli/U7EukfYG/"

$ws.Cells.Item(104, 1).Value = 103
$ws.Cells.Item(104, 2).Value = "italy"
$ws.Cells.Item(104, 3).Value = "serie-a"
$ws.Cells.Item(104, 4).Value = "2023-2024"
$ws.Cells.Item(104, 5).Value = 45234.75
$ws.Cells.Item(104, 6).Value = "Atalanta"
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = "Inter"
$ws.Cells.Item(104, 9).Value = 2
$ws.Cells.Item(104, 10).Value = 2.75
$ws.Cells.Item(104, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(104, 12).Value = 3.62
$ws.Cells.Item(104, 13).Value = "04/11/2023 17:58"
$ws.Cells.Item(104, 14).Value = 3.43
$ws.Cells.Item(104, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(104, 16).Value = 3.48
$ws.Cells.Item(104, 17).Value = "04/11/2023 17:49"
$ws.Cells.Item(104, 18).Value = 2.64
$ws.Cells.Item(104, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(104, 20).Value = 2.15
$ws.Cells.Item(104, 21).Value = "04/11/2023 17:50"
$ws.Cells.Item(104, 22).Value = "https://www.betexplorer.com/football/italy/serie-a/atalanta-inter/rcwBZEYT/"

$ws.Cells.Item(105, 1).Value = 104
$ws.Cells.Item(105, 2).Value = "italy"
$ws.Cells.Item(105, 3).Value = "serie-a"
$ws.Cells.Item(105, 4).Value = "2023-2024"
$ws.Cells.Item(105, 5).Value = 45234.86458333334
$ws.Cells.Item(105, 6).Value = "AC Milan"
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = "Udinese"
$ws.Cells.Item(105, 9).Value = 1
$ws.Cells.Item(105, 10).Value = 1.5
$ws.Cells.Item(105, 11).Value = "22/10/2023 12:02"
$ws.Cells.Item(105, 12).Value = 1.44
$ws.Cells.Item(105, 13).Value = "04/11/2023 20:44"
$ws.Cells.Item(105, 14).Value = 4.44
$ws.Cells.Item(105, 15).Value = "22/10/2023 12:02"
$ws.Cells.Item(105, 16).Value = 4.58
$ws.Cells.Item(105, 17).Value = "04/11/2023 20:44"
$ws.Cells.Item(105, 18).Value = 6.86
$ws.Cells.Item(105, 19).Value = "22/10/2023 12:02"
$ws.Cells.Item(105, 20).Value = 8.23
$ws.Cells.Item(105, 21).Value = "04/11/2023 20:44"
$ws.Cells.Item(105, 22).Value = "https://www.betexplorer.com/football/italy/serie-a/ac-milan-udinese/YuLWjG34/"
